$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 376, which shifts existing rows 376-413 down to 377-414.
$ws.Rows.Item(376).Insert()

# Populate the newly inserted row 376 with a new weekly data point.
# All fields other than Fecha (D) and Volumen (J) are copied from the
# (now shifted) row that used to occupy position 376.
$ws.Range("A376").Value2 = 4
$ws.Range("B376").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C376").Value2 = "Los Lagos"
$ws.Range("D376").Value2 = 45212
$ws.Range("E376").Value2 = 10
$ws.Range("F376").Value2 = 100112039
$ws.Range("G376").Value2 = "Ciboulette"
$ws.Range("H376").Value2 = "Sin especificar"
$ws.Range("I376").Value2 = "Primera"
$ws.Range("J376").Value2 = 240
$ws.Range("K376").Value2 = 3500
$ws.Range("L376").Value2 = 3500
$ws.Range("M376").Value2 = 3500
$ws.Range("N376").Value2 = "$/docena de atados"
$ws.Range("O376").Value2 = "Región Metropolitana"
$ws.Range("P376").Value2 = 1167
$ws.Range("Q376").Value2 = 3
$ws.Range("R376").Value2 = "Hortaliza"
